# Update the "type" (column B) and "value" (column C) columns for every
# data row (rows 2-128) according to the new naming scheme:
#   club-sports + Cross-Country-Boys  -> sports_club_boys + Cross Country
#   club-sports + Cross-Country-Girls -> sports_club_girls + Cross Country
#   club-sports + Cross-Country-Coed  -> sports_club_coed  + Cross Country
#   uil-sports  + Cross-Country-Boys  -> sports_uil_boys   + Cross Country
#   uil-sports  + Cross-Country-Girls -> sports_uil_girls  + Cross Country
#   uil-sports  + Cross-Country-Coed  -> sports_uil_coed   + Cross Country

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 2; $r -le $lastRow; $r++) {
    $typeCell = $ws.Cells.Item($r, 2)
    $valueCell = $ws.Cells.Item($r, 3)

    $typeVal = $typeCell.Value2
    $valueVal = $valueCell.Value2

    if ([string]::IsNullOrEmpty($typeVal)) { continue }

    if ($typeVal -eq "club-sports") {
        if ($valueVal -eq "Cross-Country-Boys") {
            $typeCell.Value = "sports_club_boys"
            $valueCell.Value = "Cross Country"
        } elseif ($valueVal -eq "Cross-Country-Girls") {
            $typeCell.Value = "sports_club_girls"
            $valueCell.Value = "Cross Country"
        } elseif ($valueVal -eq "Cross-Country-Coed") {
            $typeCell.Value = "sports_club_coed"
            $valueCell.Value = "Cross Country"
        }
    } elseif ($typeVal -eq "uil-sports") {
        if ($valueVal -eq "Cross-Country-Boys") {
            $typeCell.Value = "sports_uil_boys"
            $valueCell.Value = "Cross Country"
        } elseif ($valueVal -eq "Cross-Country-Girls") {
            $typeCell.Value = "sports_uil_girls"
            $valueCell.Value = "Cross Country"
        } elseif ($valueVal -eq "Cross-Country-Coed") {
            $typeCell.Value = "sports_uil_coed"
            $valueCell.Value = "Cross Country"
        }
    }
}
